$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stray _GoBack bookmark that currently sits, empty, in
#    the centered paragraph right after the title ("... összegzése").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Find the paragraph that ends with "... fájlt illetve a
#    dokumentációt." and add a brand-new paragraph right after it.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$anchor.Find.Execute("fájlt illetve a dokumentációt.", $false, $false, $false,
                      $false, $false, $true, 1, $false, "", 0) | Out-Null

$srcPara = $anchor.Paragraphs(1).Range
$srcPara.InsertParagraphAfter()

$newParaIndex = $srcPara.Paragraphs(1).Index + 1
$newPara = $d.Paragraphs($newParaIndex).Range
$newPara.Text = "A hármunk közös munkája volt az osztály diagram"

# ------------------------------------------------------------------
# 3) The first run - "A hármunk " - needs a single underline; the
#    rest of the sentence stays un-underlined.
# ------------------------------------------------------------------
$underlineRange = $d.Paragraphs($newParaIndex).Range.Duplicate
$underlineRange.Find.ClearFormatting()
$underlineRange.Find.Execute("A hármunk ", $false, $false, $false, $false,
                              $false, $true, 1, $false, "", 0) | Out-Null
$underlineRange.Font.Underline = 1

# ------------------------------------------------------------------
# 4) Re-create the _GoBack bookmark, collapsed, at the very end of
#    the text we just typed (this is where Word leaves it after the
#    last edit). A zero-length Range confuses Bookmarks.Add in this
#    host, so we bookmark a 1-char placeholder and then delete the
#    placeholder character, which leaves the bookmark collapsed in
#    the right spot.
# ------------------------------------------------------------------
$tailPara = $d.Paragraphs($newParaIndex).Range
$tailPos = $tailPara.End - 1

$placeholder = $d.Range($tailPos, $tailPos)
$placeholder.InsertAfter("X")

$placeholderRange = $d.Range($tailPos, $tailPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

$cleanupRange = $d.Range($tailPos, $tailPos + 1)
$cleanupRange.Text = ""
